# Add new header labels for the assignment columns F1:K1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "ms sql 1"
$ws.Range("G1").Value = "redis 2"
$ws.Range("H1").Value = "osgre 3"
$ws.Range("I1").Value = "mongo 4"
$ws.Range("J1").Value = "mysql 5"
$ws.Range("K1").Value = "hbase 6"

# Seed literal cycle 1..6 for the first block of students (rows 2-7)
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6

# Remaining students repeat the same cycle of assignment numbers by
# referencing the cell six rows above (shared formulas, matching Excel's
# own auto-fill behaviour down the column).
$ws.Range("E8:E71").Formula = "=E2"
$ws.Range("E72:E91").Formula = "=E66"

# Match the selection/viewport state left behind by the edit
[void]$ws.Range("E8:E91").Select()
